# 221107 upload: refresh the "datetimeFigureOut" date placeholders (master +
# every layout) from 2022-11-01 to 2022-11-07, and rewrite the auto-spin
# stop-condition bullet on the Rock Climber rules slide.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholders: Slide Master + all Custom (slide) Layouts.
# ---------------------------------------------------------------------
$d = $p.Designs.Item(1)
$master = $d.SlideMaster

$mShp = $master.Shapes.Item(3)
$mTr = $mShp.TextFrame.TextRange
$mSub = $mTr.Characters(1, $mTr.Length)
$mSub.Text = "2022-11-07"

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
  $layout = $layouts.Item($i)
  $shapes = $layout.Shapes
  for ($j = 1; $j -le $shapes.Count; $j++) {
    $shp = $shapes.Item($j)
    if ($shp.Name -like "*Date*") {
      $tr = $shp.TextFrame.TextRange
      $sub = $tr.Characters(1, $tr.Length)
      $sub.Text = "2022-11-07"
    }
  }
}

# ---------------------------------------------------------------------
# 2. Slide 2 ("Rock Climber Slot Rule") - rewrite the auto-spin bullet.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$ruleShp = $s2.Shapes.Item(1)
$ruleTr = $ruleShp.TextFrame.TextRange

$fullText = $ruleTr.Text
$startIdx = $fullText.IndexOf("Pressing spin button")
$stopIdx = $fullText.IndexOf("Win game procedure")
$oldLen = $stopIdx - $startIdx

$bullet = $ruleTr.Characters($startIdx + 1, $oldLen)
$bullet.Text = "Pressing spin button for 2 secs triggers auto spin. It will be stop with pressing stop button."

$fullText2 = $ruleTr.Text
$spaceIdx = $fullText2.IndexOf(" will be stop")
$spaceRun = $ruleTr.Characters($spaceIdx + 1, 1)
$spaceRun.Font.Size = 10

$tailIdx = $fullText2.IndexOf("with pressing stop button")
$tailLen = "with pressing stop button.".Length
$tailRun = $ruleTr.Characters($tailIdx + 1, $tailLen)
$tailRun.Font.Size = 10

Write-Host "edit.ps1 applied"
